# Auto-generated edit script applying the Ixion_Profits.xlsx diff
# to the corresponding worksheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR)
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 313.58334
$ws.Range("I11").Value = 313.58334
$ws.Range("K11").Value = 313.58334
$ws.Range("M11").Value = -173.58334

$ws.Range("H132").Value = 14494649
$ws.Range("I132").Value = 1541.2858
$ws.Range("J132").Value = 37039484
$ws.Range("K132").Value = 4623.857400000001
$ws.Range("L132").Value = 111118452
$ws.Range("M132").Value = -2093.857400000001
$ws.Range("N132").Value = -111123512

$ws.Range("H137").Value = 1688.0938
$ws.Range("I137").Value = 1332.6666
$ws.Range("J137").Value = 2754.375
$ws.Range("K137").Value = 3997.9998
$ws.Range("L137").Value = 8263.125
$ws.Range("M137").Value = -1447.9998
$ws.Range("N137").Value = -13363.125

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2839
$ws.Range("I2").Value = 1972
$ws.Range("K2").Value = 1972
$ws.Range("M2").Value = -1859

$ws.Range("H116").Value = 2839
$ws.Range("I116").Value = 1972
$ws.Range("K116").Value = 1972
$ws.Range("M116").Value = 322

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2839
$ws.Range("I3").Value = 1972
$ws.Range("K3").Value = 1972
$ws.Range("M3").Value = -1858

$ws.Range("H20").Value = 17988.4
$ws.Range("I20").Value = 1316.5714
$ws.Range("J20").Value = 32576.25
$ws.Range("K20").Value = 1316.5714
$ws.Range("L20").Value = 32576.25
$ws.Range("M20").Value = -1069.5714
$ws.Range("N20").Value = -33070.25

$ws.Range("H86").Value = 1798.4615
$ws.Range("I86").Value = 1760
$ws.Range("J86").Value = 1860
$ws.Range("K86").Value = 1760
$ws.Range("L86").Value = 1860
$ws.Range("M86").Value = -637
$ws.Range("N86").Value = -4106

$ws.Range("H89").Value = 1798.4615
$ws.Range("I89").Value = 1760
$ws.Range("J89").Value = 1860
$ws.Range("K89").Value = 8800
$ws.Range("L89").Value = 9300
$ws.Range("M89").Value = -3184
$ws.Range("N89").Value = -20532

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 6149827
$ws.Range("I31").Value = 1469.3778
$ws.Range("J31").Value = 13244086
$ws.Range("K31").Value = 1469.3778
$ws.Range("L31").Value = 13244086
$ws.Range("M31").Value = -1174.3778
$ws.Range("N31").Value = -13244676

$ws.Range("H34").Value = 6149827
$ws.Range("I34").Value = 1469.3778
$ws.Range("J34").Value = 13244086
$ws.Range("K34").Value = 1469.3778
$ws.Range("L34").Value = 13244086
$ws.Range("M34").Value = -1267.3778
$ws.Range("N34").Value = -13244490

$ws.Range("H63").Value = 41400
$ws.Range("J63").Value = 41400
$ws.Range("L63").Value = 41400
$ws.Range("N63").Value = -42772

$ws.Range("H66").Value = 41400
$ws.Range("J66").Value = 41400
$ws.Range("L66").Value = 124200
$ws.Range("N66").Value = -131064

$ws.Range("H74").Value = 19739
$ws.Range("I74").Value = 1285
$ws.Range("J74").Value = 22375.285
$ws.Range("K74").Value = 1285
$ws.Range("L74").Value = 22375.285
$ws.Range("M74").Value = -411
$ws.Range("N74").Value = -24123.285

$ws.Range("H77").Value = 19739
$ws.Range("I77").Value = 1285
$ws.Range("J77").Value = 22375.285
$ws.Range("K77").Value = 3855
$ws.Range("L77").Value = 67125.855
$ws.Range("M77").Value = 513
$ws.Range("N77").Value = -75861.855

$ws.Range("H94").Value = 5463.278
$ws.Range("I94").Value = 6458
$ws.Range("J94").Value = 4830.273
$ws.Range("K94").Value = 6458
$ws.Range("L94").Value = 4830.273
$ws.Range("M94").Value = -6007
$ws.Range("N94").Value = -5732.273

$ws.Range("H122").Value = 2559.4167
$ws.Range("I122").Value = 1745.8889
$ws.Range("J122").Value = 5000
$ws.Range("K122").Value = 5237.6667
$ws.Range("L122").Value = 15000
$ws.Range("M122").Value = -2787.6667
$ws.Range("N122").Value = -19900

$ws.Range("H132").Value = 6454030
$ws.Range("I132").Value = 9525586
$ws.Range("J132").Value = 3763
$ws.Range("K132").Value = 28576758
$ws.Range("L132").Value = 11289
$ws.Range("M132").Value = -28574228
$ws.Range("N132").Value = -16349

$ws.Range("H134").Value = 8549322
$ws.Range("I134").Value = 11907388
$ws.Range("J134").Value = 1518.909
$ws.Range("K134").Value = 35722164
$ws.Range("L134").Value = 4556.727000000001
$ws.Range("M134").Value = -35719629
$ws.Range("N134").Value = -9626.727000000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H33").Value = 4352242.5
$ws.Range("I33").Value = 5882422
$ws.Range("J33").Value = 16734.5
$ws.Range("K33").Value = 35294532
$ws.Range("L33").Value = 100407
$ws.Range("M33").Value = -35294249
$ws.Range("N33").Value = -100973

$ws.Range("H113").Value = 1818755.2
$ws.Range("I113").Value = 1724718.1
$ws.Range("K113").Value = 5174154.300000001
$ws.Range("M113").Value = -5171984.300000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H14").Value = 3671200
$ws.Range("I14").Value = 5501800
$ws.Range("K14").Value = 5501800
$ws.Range("M14").Value = -5501632

$ws.Range("H109").Value = 0
$ws.Range("J109").Value = 0
$ws.Range("L109").Value = 0
$ws.Range("N109").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 55556404
$ws.Range("I93").Value = 956.375
$ws.Range("J93").Value = 500000000
$ws.Range("K93").Value = 956.375
$ws.Range("L93").Value = 500000000
$ws.Range("M93").Value = 291.625
$ws.Range("N93").Value = -500002496

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 1482.8235
$ws.Range("I122").Value = 1325.5
$ws.Range("J122").Value = 4000
$ws.Range("K122").Value = 3976.5
$ws.Range("L122").Value = 12000
$ws.Range("M122").Value = -1526.5
$ws.Range("N122").Value = -16900

$ws.Range("H132").Value = 1281.6444
$ws.Range("I132").Value = 851.96155
$ws.Range("J132").Value = 1869.6316
$ws.Range("K132").Value = 2555.88465
$ws.Range("L132").Value = 5608.8948
$ws.Range("M132").Value = -25.88464999999997
$ws.Range("N132").Value = -10668.8948

$ws.Range("H136").Value = 19611820
$ws.Range("I136").Value = 5634
$ws.Range("J136").Value = 41668780
$ws.Range("K136").Value = 16902
$ws.Range("L136").Value = 125006340
$ws.Range("M136").Value = -14352
$ws.Range("N136").Value = -125011440
